# Generate Report for Handoff
# Rename the two localization entities and refresh their handoff metadata:
#   013259bb-e739-43cb-b1bc-3d564cffbc96  ->  3e063fe6-440a-40e3-b126-f6157dce417f
#   7bd24413-cf0e-405e-a79a-ae2cd1cdab88  ->  ffffe006dce6-1480-4b67-aa8d-a430f64d017f
# The second entry's handoff xliff now matches the first entry's (content duplicate),
# so "Content Duplicate" flips to True and the Latest Handoff File / Datetime columns
# for row 3 are refreshed to the same values as row 2.

$wb = $excel.ActiveWorkbook

$oldId1 = "013259bb-e739-43cb-b1bc-3d564cffbc96"
$newId1 = "3e063fe6-440a-40e3-b126-f6157dce417f"
$oldId2 = "7bd24413-cf0e-405e-a79a-ae2cd1cdab88"
$newId2 = "ffffe006dce6-1480-4b67-aa8d-a430f64d017f"

$newZhXlf = "$newId1.529674f7dbc67588d900055118d63c0bad1e9569.zh-cn.xlf"
$newDeXlf = "$newId1.529674f7dbc67588d900055118d63c0bad1e9569.de-de.xlf"

$newOverviewDate = "2017-01-03 05:36:38"
$newZhHandoffDate = "2017-01-03 05:36:27"
$newDeHandoffDate = "2017-01-03 05:36:38"

# ---- Sheet "Overview" ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newId1.md"
$wsOverview.Range("B2").Value = "e2e\$newId1.md"
$wsOverview.Range("G2").Value = $newOverviewDate
$wsOverview.Range("A3").Value = "$newId2.md"
$wsOverview.Range("B3").Value = "e2e\$newId2.md"
$wsOverview.Range("G3").Value = $newOverviewDate

$wsOverview.Hyperlinks.Item(1).TextToDisplay = "e2e\$newId1.md"
$wsOverview.Hyperlinks.Item(2).TextToDisplay = "e2e\$newId2.md"

# ---- Sheet "zh-cn" ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = "$newId1.md"
$wsZh.Range("G2").Value = $newZhXlf
$wsZh.Range("H2").Value = $newZhHandoffDate
$wsZh.Range("A3").Value = "$newId2.md"
$wsZh.Range("F3").Value = "'True"
$wsZh.Range("G3").Value = $newZhXlf
$wsZh.Range("H3").Value = $newZhHandoffDate

$wsZh.Hyperlinks.Item(1).TextToDisplay = "$newId1.md"
$wsZh.Hyperlinks.Item(2).TextToDisplay = "$newId2.md"

# ---- Sheet "de-de" ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = "$newId1.md"
$wsDe.Range("G2").Value = $newDeXlf
$wsDe.Range("H2").Value = $newDeHandoffDate
$wsDe.Range("A3").Value = "$newId2.md"
$wsDe.Range("F3").Value = "'True"
$wsDe.Range("G3").Value = $newDeXlf
$wsDe.Range("H3").Value = $newDeHandoffDate

$wsDe.Hyperlinks.Item(1).TextToDisplay = "$newId1.md"
$wsDe.Hyperlinks.Item(2).TextToDisplay = "$newId2.md"
